# "Rows transposed to columns": the year header row (B1:K1) and the
# "Domestic travel" data row (B2:K2) become two columns (year, Domestic
# travel) running down the sheet. Only 2008-2015 (old D:K) carry over;
# 1998/2003 (old B/C) have no counterpart in the new layout and are
# dropped, which also shortens the sheet by two trailing blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$years  = @(2008, 2009, 2010, 2011, 2012, 2013, 2014, 2015)
$values = @(110.22, 144.202, 116.578, 133.87899999999999, 117.465, 132.608, 129.017, 128.5)

# Row 2 used to be the "Domestic travel" row with a custom 20.1pt height;
# it becomes a plain data row now, so restore the default row height
# before touching its formatting/content.
$ws.Rows.Item(2).AutoFit()

# New column B header ("      Domestic travel", formerly A2) picks up the
# same style A1 uses (s="2" = centered header style).
$ws.Range("B1").Value = "      Domestic travel"
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write the transposed (year, value) pairs down columns A/B starting at
# row 2, clearing away the placeholder/label style as we go.
for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Range("A" + $row + ":B" + $row).ClearFormats()
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Drop the old year/value columns C:K from rows 1 and 2.
$ws.Range("C1:K2").ClearContents()

# The sheet now only needs 28 rows instead of 30 - remove the two extra
# trailing blank rows.
$ws.Rows.Item(30).Delete()
$ws.Rows.Item(29).Delete()

# Column B gets its own width now that it holds the "Domestic travel"
# values/header (mirrors the bestFit width column A already has).
$ws.Columns.Item(2).ColumnWidth = 16.8

# Match the post-edit selection (the new data block, A2:XFD3).
$ws.Range("A2:XFD3").Select()
